$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the existing "EstMusical" lookup table 6 rows down     ---
# (R1:S5 -> R7:S11) to make room for the new "Visualizacao" table.
$ws.Range("R1:S5").Cut($ws.Range("R7:S11"))
$ws.Range("S7").ClearContents()

# --- New "Visualizacao" table headers (row 2) ---
$ws.Range("P2").Value = "dataLançamento"
$ws.Range("Q2").Value = "quantMin"
$ws.Range("R2").Value = "localizacao"
$ws.Range("S2").Value = "Visualização"

# --- New "Visualizacao" table data (rows 3-4) ---
$ws.Range("P3").Value = 44413
$ws.Range("P4").Value = 44414

$ws.Range("Q4").Value = "50min"
$ws.Range("Q3").Value = "30min"
$ws.Range("R3").Value = "Mato Grosso"
$ws.Range("R4").Value = "Rio de Janeiro"
$ws.Range("S3").Value = "Desativo"
$ws.Range("S4").Value = "Ativo"

# Apply the date number format to P3 then copy it (format-only) to P4 so
# both cells share a single cellXf record, matching a real Excel fill/paste.
$ws.Range("P3").NumberFormat = "mm-dd-yy"
$ws.Range("P3").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View/selection tweaks ---
$ws.Range("S4").Select()
